# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.318.60"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.276.51"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.99"
$ws.Range("E5").Value = "  -5.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.68"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -0.48%  "

$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("E9").Value = "  -1.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.32"
$ws.Range("E10").Value = "  -3.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0896"
$ws.Range("E11").Value = "  -1.17%  "

$ws.Range("E12").Value = "  -2.20%  "

$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.962"
$ws.Range("E14").Value = "  -0.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.06"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.621.97"
$ws.Range("E16").Value = "  -0.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.272.26"
$ws.Range("E17").Value = "  -0.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.529.17"
$ws.Range("E18").Value = "  +0.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.22"
$ws.Range("E19").Value = "  -2.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000103"
$ws.Range("E20").Value = "  -1.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.24"
$ws.Range("E21").Value = "  +2.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.53"
$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.36"
$ws.Range("E23").Value = "  -7.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.35"
$ws.Range("E24").Value = "  -2.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.16"
$ws.Range("E25").Value = "  -2.92%  "

$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.59"
$ws.Range("E27").Value = "  -2.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.34"
$ws.Range("E28").Value = "  +0.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.98"
$ws.Range("E29").Value = "  +13.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.03"
$ws.Range("E30").Value = "  -1.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.54"
$ws.Range("E31").Value = "  -5.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.83"
$ws.Range("E32").Value = "  -0.24%  "

$ws.Range("E33").Value = "  -3.72%  "

$ws.Range("E34").Value = "  -3.55%  "

$ws.Range("E35").Value = "  +1.85%  "

$ws.Range("E36").Value = "  -3.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.46"
$ws.Range("E37").Value = "  -3.07%  "

$ws.Range("E38").Value = "  -2.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.60"
$ws.Range("E39").Value = "  -3.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.69"
$ws.Range("E40").Value = "  -1.91%  "

$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.01"
$ws.Range("E42").Value = "  +8.75%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("E43").Value = "  +0.56%  "

$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "68.57"
$ws.Range("E44").Value = "  -1.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.224"
$ws.Range("E45").Value = "  -0.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.91"
$ws.Range("E46").Value = "  -3.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.703.43"
$ws.Range("E47").Value = "  +7.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "109.44"
$ws.Range("E48").Value = "  -2.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.07"
$ws.Range("E49").Value = "  -5.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.60"
$ws.Range("E50").Value = "  -3.73%  "

$ws.Range("E51").Value = "  -2.26%  "
